$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 82.38
$ws.Range("I15").Value = 82.38
$ws.Range("K15").Value = 247.14
$ws.Range("M15").Value = -78.13999999999999
$ws.Range("H17").Value = 179695.95
$ws.Range("J17").Value = 179695.95
$ws.Range("L17").Value = 539087.8500000001
$ws.Range("N17").Value = -539423.8500000001
$ws.Range("H116").Value = 2760.5652
$ws.Range("I116").Value = 2527.389
$ws.Range("J116").Value = 3600
$ws.Range("K116").Value = 2527.389
$ws.Range("L116").Value = 3600
$ws.Range("M116").Value = 914.6109999999999
$ws.Range("N116").Value = -10484
$ws.Range("H132").Value = 1192.375
$ws.Range("J132").Value = 680.6
$ws.Range("L132").Value = 2041.8
$ws.Range("N132").Value = -7101.8
$ws.Range("H135").Value = 715.1719000000001
$ws.Range("I135").Value = 393.78183
$ws.Range("K135").Value = 3544.03647
$ws.Range("M135").Value = -1009.03647
$ws.Range("H137").Value = 796.8946999999999
$ws.Range("I137").Value = 653.26666
$ws.Range("J137").Value = 890.56525
$ws.Range("K137").Value = 1959.79998
$ws.Range("L137").Value = 2671.69575
$ws.Range("M137").Value = 590.20002
$ws.Range("N137").Value = -7771.69575
$ws.Range("H138").Value = 1830.49
$ws.Range("I138").Value = 833.42224
$ws.Range("J138").Value = 2646.2727
$ws.Range("K138").Value = 2500.26672
$ws.Range("L138").Value = 7938.8181
$ws.Range("M138").Value = 2639.73328
$ws.Range("N138").Value = -18218.8181
$ws.Range("H141").Value = 2224.2979
$ws.Range("I141").Value = 895.6923
$ws.Range("K141").Value = 2687.0769
$ws.Range("M141").Value = 2492.9231

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20578.031
$ws.Range("I32").Value = 19702.805
$ws.Range("J32").Value = 24322.055
$ws.Range("K32").Value = 19702.805
$ws.Range("L32").Value = 24322.055
$ws.Range("M32").Value = -19415.805
$ws.Range("N32").Value = -24896.055
$ws.Range("H61").Value = 888.09753
$ws.Range("I61").Value = 753.82355
$ws.Range("J61").Value = 1540.2858
$ws.Range("K61").Value = 753.82355
$ws.Range("L61").Value = 1540.2858
$ws.Range("M61").Value = -541.82355
$ws.Range("N61").Value = -1964.2858
$ws.Range("H74").Value = 869.125
$ws.Range("I74").Value = 870.44684
$ws.Range("J74").Value = 862.2222
$ws.Range("K74").Value = 870.44684
$ws.Range("L74").Value = 862.2222
$ws.Range("M74").Value = 3.553160000000048
$ws.Range("N74").Value = -2610.2222
$ws.Range("H77").Value = 869.125
$ws.Range("I77").Value = 870.44684
$ws.Range("J77").Value = 862.2222
$ws.Range("K77").Value = 4352.2342
$ws.Range("L77").Value = 4311.111
$ws.Range("M77").Value = 15.76580000000013
$ws.Range("N77").Value = -13047.111
$ws.Range("H102").Value = 3803.3333
$ws.Range("I102").Value = 3764
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 3764
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -2142
$ws.Range("N102").Value = -7244
$ws.Range("H122").Value = 1253
$ws.Range("I122").Value = 902.8333
$ws.Range("J122").Value = 2303.5
$ws.Range("K122").Value = 2708.4999
$ws.Range("L122").Value = 6910.5
$ws.Range("M122").Value = -258.4998999999998
$ws.Range("N122").Value = -11810.5
$ws.Range("H123").Value = 43320
$ws.Range("J123").Value = 43320
$ws.Range("L123").Value = 43320
$ws.Range("N123").Value = -53120
$ws.Range("H132").Value = 1098.7441
$ws.Range("I132").Value = 980.9091
$ws.Range("J132").Value = 1487.6
$ws.Range("K132").Value = 2942.7273
$ws.Range("L132").Value = 4462.799999999999
$ws.Range("M132").Value = -412.7273
$ws.Range("N132").Value = -9522.799999999999
$ws.Range("H136").Value = 888.09753
$ws.Range("I136").Value = 753.82355
$ws.Range("J136").Value = 1540.2858
$ws.Range("K136").Value = 2261.47065
$ws.Range("L136").Value = 4620.857400000001
$ws.Range("M136").Value = 288.5293500000002
$ws.Range("N136").Value = -9720.857400000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2102.95
$ws.Range("I86").Value = 1998.8529
$ws.Range("J86").Value = 2692.8333
$ws.Range("K86").Value = 1998.8529
$ws.Range("L86").Value = 2692.8333
$ws.Range("M86").Value = -875.8529000000001
$ws.Range("N86").Value = -4938.8333
$ws.Range("H89").Value = 2102.95
$ws.Range("I89").Value = 1998.8529
$ws.Range("J89").Value = 2692.8333
$ws.Range("K89").Value = 9994.264500000001
$ws.Range("L89").Value = 13464.1665
$ws.Range("M89").Value = -4378.264500000001
$ws.Range("N89").Value = -24696.1665
$ws.Range("H99").Value = 27779242
$ws.Range("I99").Value = 33334712
$ws.Range("K99").Value = 33334712
$ws.Range("M99").Value = -33333214
$ws.Range("H105").Value = 4339.5
$ws.Range("I105").Value = 3630.9048
$ws.Range("J105").Value = 6465.2856
$ws.Range("K105").Value = 3630.9048
$ws.Range("L105").Value = 6465.2856
$ws.Range("M105").Value = -1883.9048
$ws.Range("N105").Value = -9959.285599999999
$ws.Range("H134").Value = 14314.705
$ws.Range("I134").Value = 1312.209
$ws.Range("J134").Value = 93511.73
$ws.Range("K134").Value = 3936.627
$ws.Range("L134").Value = 280535.19
$ws.Range("M134").Value = -1401.627
$ws.Range("N134").Value = -285605.19

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 141.46666
$ws.Range("I7").Value = 156.45454
$ws.Range("J7").Value = 100.25
$ws.Range("K7").Value = 156.45454
$ws.Range("L7").Value = 100.25
$ws.Range("M7").Value = -43.45454000000001
$ws.Range("N7").Value = -326.25
$ws.Range("H31").Value = 2978
$ws.Range("I31").Value = 2805.4285
$ws.Range("J31").Value = 3783.3333
$ws.Range("K31").Value = 2805.4285
$ws.Range("L31").Value = 3783.3333
$ws.Range("M31").Value = -2510.4285
$ws.Range("N31").Value = -4373.3333
$ws.Range("H34").Value = 2978
$ws.Range("I34").Value = 2805.4285
$ws.Range("J34").Value = 3783.3333
$ws.Range("K34").Value = 2805.4285
$ws.Range("L34").Value = 3783.3333
$ws.Range("M34").Value = -2603.4285
$ws.Range("N34").Value = -4187.3333
$ws.Range("H58").Value = 3997.879
$ws.Range("I58").Value = 856.5263
$ws.Range("J58").Value = 8261.143
$ws.Range("K58").Value = 856.5263
$ws.Range("L58").Value = 8261.143
$ws.Range("M58").Value = -653.5263
$ws.Range("N58").Value = -8667.143
$ws.Range("H94").Value = 2981.9473
$ws.Range("I94").Value = 3484.8572
$ws.Range("J94").Value = 2868.3872
$ws.Range("K94").Value = 3484.8572
$ws.Range("L94").Value = 2868.3872
$ws.Range("M94").Value = -3033.8572
$ws.Range("N94").Value = -3770.3872
$ws.Range("H132").Value = 1511.8025
$ws.Range("I132").Value = 933.8
$ws.Range("J132").Value = 2444.0645
$ws.Range("K132").Value = 2801.4
$ws.Range("L132").Value = 7332.193499999999
$ws.Range("M132").Value = -271.3999999999996
$ws.Range("N132").Value = -12392.1935
$ws.Range("H134").Value = 1288.3793
$ws.Range("I134").Value = 1264.4615
$ws.Range("J134").Value = 1337.4736
$ws.Range("K134").Value = 3793.3845
$ws.Range("L134").Value = 4012.4208
$ws.Range("M134").Value = -1258.3845
$ws.Range("N134").Value = -9082.4208
$ws.Range("H136").Value = 3997.879
$ws.Range("I136").Value = 856.5263
$ws.Range("J136").Value = 8261.143
$ws.Range("K136").Value = 2569.5789
$ws.Range("L136").Value = 24783.429
$ws.Range("M136").Value = -19.57889999999998
$ws.Range("N136").Value = -29883.429

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 8000
$ws.Range("J93").Value = 8000
$ws.Range("L93").Value = 24000
$ws.Range("N93").Value = -27744
$ws.Range("H122").Value = 494.10938
$ws.Range("I122").Value = 269.14285
$ws.Range("J122").Value = 557.1
$ws.Range("K122").Value = 2422.28565
$ws.Range("L122").Value = 5013.900000000001
$ws.Range("M122").Value = 27.71434999999974
$ws.Range("N122").Value = -9913.900000000001
$ws.Range("H131").Value = 30843.045
$ws.Range("J131").Value = 18176.432
$ws.Range("L131").Value = 54529.296
$ws.Range("N131").Value = -64609.296

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4829.9067
$ws.Range("I70").Value = 4806.2
$ws.Range("J70").Value = 4884.615
$ws.Range("K70").Value = 4806.2
$ws.Range("L70").Value = 4884.615
$ws.Range("M70").Value = -4536.2
$ws.Range("N70").Value = -5424.615
$ws.Range("H73").Value = 4829.9067
$ws.Range("I73").Value = 4806.2
$ws.Range("J73").Value = 4884.615
$ws.Range("K73").Value = 4806.2
$ws.Range("L73").Value = 4884.615
$ws.Range("M73").Value = -3870.2
$ws.Range("N73").Value = -6756.615
$ws.Range("H97").Value = 3140
$ws.Range("I97").Value = 3140
$ws.Range("K97").Value = 3140
$ws.Range("M97").Value = -2644
$ws.Range("H122").Value = 25965392
$ws.Range("I122").Value = 26608588
$ws.Range("J122").Value = 25000600
$ws.Range("K122").Value = 79825764
$ws.Range("L122").Value = 75001800
$ws.Range("M122").Value = -79823314
$ws.Range("N122").Value = -75006700
$ws.Range("H132").Value = 1785.4286
$ws.Range("I132").Value = 1692.8611
$ws.Range("J132").Value = 1908.8518
$ws.Range("K132").Value = 5078.5833
$ws.Range("L132").Value = 5726.555399999999
$ws.Range("M132").Value = -2548.5833
$ws.Range("N132").Value = -10786.5554

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1852.4688
$ws.Range("I46").Value = 1209.3334
$ws.Range("J46").Value = 2419.9412
$ws.Range("K46").Value = 1209.3334
$ws.Range("L46").Value = 2419.9412
$ws.Range("M46").Value = -1021.3334
$ws.Range("N46").Value = -2795.9412
$ws.Range("H100").Value = 15875630
$ws.Range("I100").Value = 22224942
$ws.Range("J100").Value = 2350
$ws.Range("K100").Value = 22224942
$ws.Range("L100").Value = 2350
$ws.Range("M100").Value = -22224401
$ws.Range("N100").Value = -3432
$ws.Range("H122").Value = 13645.777
$ws.Range("I122").Value = 16687.428
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 50062.284
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -47612.284
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 2068.9412
$ws.Range("I132").Value = 1867.925
$ws.Range("J132").Value = 2799.9092
$ws.Range("K132").Value = 5603.775
$ws.Range("L132").Value = 8399.7276
$ws.Range("M132").Value = -3073.775
$ws.Range("N132").Value = -13459.7276

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 862.9167
$ws.Range("I122").Value = 810.5172
$ws.Range("J122").Value = 1080
$ws.Range("K122").Value = 2431.5516
$ws.Range("L122").Value = 3240
$ws.Range("M122").Value = 18.44840000000022
$ws.Range("N122").Value = -8140
$ws.Range("H132").Value = 1130.8286
$ws.Range("I132").Value = 1040
$ws.Range("J132").Value = 1284.5385
$ws.Range("K132").Value = 3120
$ws.Range("L132").Value = 3853.6155
$ws.Range("M132").Value = -590
$ws.Range("N132").Value = -8913.6155
$ws.Range("H136").Value = 1371.5652
$ws.Range("I136").Value = 1532.1765
$ws.Range("J136").Value = 916.5
$ws.Range("K136").Value = 4596.529500000001
$ws.Range("L136").Value = 2749.5
$ws.Range("M136").Value = -2046.529500000001
$ws.Range("N136").Value = -7849.5

Write-Output "Applied 285 cell updates across 8 sheets"